# Add data for 2024-10-25
# Updates 2024 (column K) totals across Citywide Totals, By Neighborhood,
# and individual neighborhood sheets for the crimes recorded on 2024-10-25.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6617   # Aggravated Assault: 6601 -> 6617
$ws.Range("K3").Value = 6842   # Aggravated Battery: 6821 -> 6842
$ws.Range("K4").Value = 1417   # Criminal Sexual Assault: 1414 -> 1417
$ws.Range("K6").Value = 7524   # Robbery: 7496 -> 7524
$ws.Range("K7").Value = 22895   # Total: 22827 -> 22895

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 412   # Aggravated Assault: 410 -> 412
$ws.Range("K7").Value = 1499   # Total: 1497 -> 1499

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K6").Value = 87   # Robbery: 86 -> 87
$ws.Range("K7").Value = 376   # Total: 375 -> 376

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 222   # Aggravated Assault: 221 -> 222
$ws.Range("K6").Value = 231   # Robbery: 230 -> 231
$ws.Range("K7").Value = 777   # Total: 775 -> 777

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K6").Value = 97   # Robbery: 96 -> 97
$ws.Range("K7").Value = 385   # Total: 384 -> 385

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K4").Value = 84   # Archer Heights: 83 -> 84
$ws.Range("K7").Value = 690   # Auburn Gresham: 684 -> 690
$ws.Range("K8").Value = 1499   # Austin: 1497 -> 1499
$ws.Range("K11").Value = 422   # Belmont Cragin: 420 -> 422
$ws.Range("K17").Value = 44   # Burnside: 43 -> 44
$ws.Range("K18").Value = 153   # Calumet Heights: 151 -> 153
$ws.Range("K19").Value = 668   # Chatham: 665 -> 668
$ws.Range("K20").Value = 552   # Chicago Lawn: 548 -> 552
$ws.Range("K27").Value = 213   # Edgewater: 209 -> 213
$ws.Range("K29").Value = 1242   # Englewood: 1238 -> 1242
$ws.Range("K31").Value = 253   # Gage Park: 252 -> 253
$ws.Range("K37").Value = 777   # Grand Crossing: 775 -> 777
$ws.Range("K42").Value = 845   # Humboldt Park: 842 -> 845
$ws.Range("K43").Value = 187   # Hyde Park: 186 -> 187
$ws.Range("K48").Value = 290   # Lake View: 288 -> 290
$ws.Range("K52").Value = 608   # Little Village: 606 -> 608
$ws.Range("K54").Value = 453   # Loop: 451 -> 453
$ws.Range("K63").Value = 63   # NO NEIGHBORHOOD DATA: 60 -> 63
$ws.Range("K67").Value = 892   # North Lawndale: 891 -> 892
$ws.Range("K76").Value = 308   # River North: 306 -> 308
$ws.Range("K77").Value = 156   # Riverdale: 154 -> 156
$ws.Range("K79").Value = 571   # Roseland: 568 -> 571
$ws.Range("K85").Value = 1058   # South Shore: 1053 -> 1058
$ws.Range("K86").Value = 141   # Streeterville: 140 -> 141
$ws.Range("K90").Value = 217   # Washington Heights: 216 -> 217
$ws.Range("K91").Value = 272   # Washington Park: 270 -> 272
$ws.Range("K93").Value = 87   # West Lawn: 85 -> 87
$ws.Range("K94").Value = 306   # West Loop: 305 -> 306
$ws.Range("K95").Value = 376   # West Pullman: 375 -> 376
$ws.Range("K98").Value = 114   # Wicker Park: 113 -> 114
$ws.Range("K99").Value = 385   # Woodlawn: 384 -> 385
$ws.Range("K101").Value = 22895   # Total: 22827 -> 22895

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K6").Value = 88   # Robbery: 87 -> 88
$ws.Range("K7").Value = 253   # Total: 252 -> 253

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K6").Value = 254   # Robbery: 253 -> 254
$ws.Range("K7").Value = 892   # Total: 891 -> 892

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 244   # Robbery: 242 -> 244
$ws.Range("K7").Value = 453   # Total: 451 -> 453

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 351   # Aggravated Assault: 350 -> 351
$ws.Range("K3").Value = 439   # Aggravated Battery: 436 -> 439
$ws.Range("K7").Value = 1242   # Total: 1238 -> 1242

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 138   # Robbery: 136 -> 138
$ws.Range("K7").Value = 290   # Total: 288 -> 290

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 200   # Aggravated Battery: 198 -> 200
$ws.Range("K6").Value = 221   # Robbery: 220 -> 221
$ws.Range("K7").Value = 668   # Total: 665 -> 668

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 71   # Aggravated Assault: 70 -> 71
$ws.Range("K3").Value = 59   # Aggravated Battery: 58 -> 59
$ws.Range("K7").Value = 308   # Total: 306 -> 308

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 226   # Aggravated Assault: 224 -> 226
$ws.Range("K6").Value = 313   # Robbery: 312 -> 313
$ws.Range("K7").Value = 845   # Total: 842 -> 845

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K3").Value = 129   # Aggravated Battery: 128 -> 129
$ws.Range("K6").Value = 56   # Robbery: 55 -> 56
$ws.Range("K7").Value = 272   # Total: 270 -> 272

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 194   # Aggravated Assault: 191 -> 194
$ws.Range("K7").Value = 571   # Total: 568 -> 571

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 191   # Aggravated Assault: 190 -> 191
$ws.Range("K3").Value = 177   # Aggravated Battery: 175 -> 177
$ws.Range("K6").Value = 151   # Robbery: 150 -> 151
$ws.Range("K7").Value = 552   # Total: 548 -> 552

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 43   # Aggravated Assault: 42 -> 43
$ws.Range("K3").Value = 50   # Aggravated Battery: 49 -> 50
$ws.Range("K7").Value = 153   # Total: 151 -> 153

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K6").Value = 13   # Robbery: 12 -> 13
$ws.Range("K7").Value = 44   # Total: 43 -> 44

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K2").Value = 27   # Aggravated Assault: 26 -> 27
$ws.Range("K6").Value = 36   # Robbery: 35 -> 36
$ws.Range("K7").Value = 87   # Total: 85 -> 87

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 225   # Aggravated Assault: 224 -> 225
$ws.Range("K3").Value = 223   # Aggravated Battery: 221 -> 223
$ws.Range("K6").Value = 190   # Robbery: 187 -> 190
$ws.Range("K7").Value = 690   # Total: 684 -> 690

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 139   # Robbery: 138 -> 139
$ws.Range("K7").Value = 306   # Total: 305 -> 306

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("K3").Value = 23   # Aggravated Battery: 22 -> 23
$ws.Range("K7").Value = 114   # Total: 113 -> 114

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 145   # Aggravated Assault: 144 -> 145
$ws.Range("K3").Value = 109   # Aggravated Battery: 108 -> 109
$ws.Range("K7").Value = 422   # Total: 420 -> 422

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 78   # Robbery: 74 -> 78
$ws.Range("K7").Value = 213   # Total: 209 -> 213

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K6").Value = 34   # Robbery: 33 -> 34
$ws.Range("K7").Value = 141   # Total: 140 -> 141

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 80   # Aggravated Assault: 79 -> 80
$ws.Range("K7").Value = 217   # Total: 216 -> 217

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K3").Value = 52   # Aggravated Battery: 51 -> 52
$ws.Range("K7").Value = 187   # Total: 186 -> 187

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 369   # Aggravated Battery: 366 -> 369
$ws.Range("K6").Value = 258   # Robbery: 256 -> 258
$ws.Range("K7").Value = 1058   # Total: 1053 -> 1058

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K3").Value = 61   # Aggravated Battery: 59 -> 61
$ws.Range("K7").Value = 156   # Total: 154 -> 156

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 163   # Aggravated Assault: 162 -> 163
$ws.Range("K4").Value = 33   # Criminal Sexual Assault: 32 -> 33
$ws.Range("K7").Value = 608   # Total: 606 -> 608

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K6").Value = 31   # Robbery: 30 -> 31
$ws.Range("K7").Value = 84   # Total: 83 -> 84
